# Add new columns I ("I0") and J ("IF") to Sheet1, mirroring the header
# style already used for column H, and populate rows 2-67 with the
# corresponding values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values for I2:I67 and J2:J67 (row 2 first .. row 67 last)
$iVals = @(8,9,9,9,8,9,9,10,9,9,9,9,9,9,9,9,9,8,8,9,9,9,9,8,8,9,9,9,8,10,8,9,9,9,9,8,8,8,8,9,7,7,6,7,12,9,7,8,8,6,6,8,4,6,5,6,6,8,7,1,5,6,9,3,6,6)
$jVals = @(8,9,9,9,9,9,9,10,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,8,8,9,9,9,9,10,8,9,9,9,9,8,8,8,9,9,7,7,6,8,14,9,7,8,8,7,6,8,4,7,6,7,8,9,7,1,6,7,9,4,6,6)

for ($r = 2; $r -le 67; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iVals[$idx]
    $ws.Cells.Item($r, 10).Value = $jVals[$idx]
}

# Copy the existing header formatting (bold, centered, bordered) from H1
# onto the two new header cells, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
